$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric need an explicit Text number format
# so Excel stores them as literal text instead of converting them to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "51.656.51"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "2.900.84"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "354.84"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "108.82"
$ws.Range("E6").Value = "  -3.22%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "38.77"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.137"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.0866"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "19.40"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "3.377.56"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "2.921.61"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "0.973"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("D18").Value = "51.717.93"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "13.86"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "70.31"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "267.10"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "0.185"
$ws.Range("E26").Value = "  +11.69%  "
$ws.Range("D27").Value = "7.64"
$ws.Range("E27").Value = "  +18.01%  "
$ws.Range("D28").Value = "26.74"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.105"
$ws.Range("E30").Value = "  +9.32%  "
$ws.Range("D31").Value = "10.44"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").Value = "37.02"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("D34").Value = "6.02"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").Value = "52.15"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").Value = "0.0440"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("D39").Value = "18.14"
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("D40").Value = "1.99"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("D41").Value = "2.70"
$ws.Range("E41").Value = "  -5.19%  "
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").Value = "22.66"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "118.93"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "2.46"
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("D48").Value = "2.121.04"
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("D49").Value = "0.250"
$ws.Range("E49").Value = "  -5.25%  "
$ws.Range("D50").Value = "0.0341"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").Value = "9.06"
$ws.Range("E51").Value = "  -0.19%  "
